$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.551.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.467.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.79%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.31%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.843.81"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.90"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.466.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.776"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.530.83"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.24"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.15"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.34"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.29"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0759"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.37%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.73%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.67%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.947.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.74"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.702.81"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.24%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.22"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.45%  "
